$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers)
$wb.Worksheets.Item(1).Name = "GNG_TO-1650291276019091"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912780746145"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912780756214"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912781299708"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291278223159"

# Sheet 1 (GNG): update B2:B5
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912759695675.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912759850729.csv"
$ws1.Range("B4").Value = "go_stims-1650291275987072.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912760167515.csv"

# Sheet 2 (NB): update B2:B10
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_9-16502912768953314.csv"
$ws2.Range("B3").Value = "ZB-match_2-16502912762733574.csv"
$ws2.Range("B4").Value = "OB-16502912771418693.csv"
$ws2.Range("B5").Value = "OB-1650291277588088.csv"
$ws2.Range("B6").Value = "TB-1650291277728235.csv"
$ws2.Range("B7").Value = "OB-16502912776175034.csv"
$ws2.Range("B8").Value = "TB-16502912780572512.csv"
$ws2.Range("B9").Value = "TB-16502912779117632.csv"
$ws2.Range("B10").Value = "ZB-match_5-1650291276579719.csv"

# Sheet 4 (TOL): update B2:B7
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912780902958.csv"
$ws4.Range("B3").Value = "ZM_stims-1650291278077617.csv"
$ws4.Range("B4").Value = "MM_stims-16502912781132019.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912780912523.csv"
$ws4.Range("B6").Value = "MM_stims-16502912781286247.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912781132019.csv"

# Sheet 5 (vSAT): update B2:B5
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650291278176647.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912781596005.csv"
$ws5.Range("B4").Value = "SAT_stims-16502912781339705.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912782073038.csv"
